# Config.xlsx update: "All CRUD and Updated data sets"
#
# - Row 2's TestDataFile (F2) is updated to point at the new data file
#   "LeaveDeductions.xlsx" instead of the old "CRUDOperations.xlsx".
# - The "RunMode" cell for that same row (B2) is given a date style
#   (DD/MM/YY) instead of the default General format.
# - The selected/active cell is moved to F2 (the cell that was edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Point the first data row at the new test-data workbook.
$ws.Range("F2").Value = "LeaveDeductions.xlsx"

# Apply a date number format to B2 (was General before).
$ws.Range("B2").NumberFormat = "DD/MM/YY"

# Leave the cursor on the cell that was changed.
$null = $ws.Range("F2").Select()
